{"js": "// Corecion TP4 (Gestion Linea Base)\n// 1) Insert a missing space after \"...fechas de ex\u00e1menes\" (before the\n//    existing spacing-adjusted run that precedes \"parciales\").\n// 2) Fix grammar in \"...terminar o bien realizar los \u00faltimos cambios\n//    sobre para Trabajos...\" -> \"...terminen de  realizar los \u00faltimos\n//    cambios sobre los Trabajos...\".\n// 3) Fix \"hallan faltado\" -> \"hayan faltado\".\n\nconst body = context.document.body;\n\n// --- Change 1: add a space after \"... fechas de ex\u00e1menes\" ---\nconst examenesHits = body.search(\"h\u00e1bil posterior a las fechas de ex\u00e1menes\", { matchCase: true });\nexamenesHits.load(\"text\");\nawait context.sync();\nif (examenesHits.items.length > 0) {\n  examenesHits.items[0].insertText(\" \", \"After\");\n  await context.sync();\n}\n\n// --- Change 2: rewrite the \"Se aprovechar\u00e1n...\" sentence ---\nconst aprovecharanHits = body.search(\n  \"Se aprovechar\u00e1n los d\u00edas h\u00e1biles para que los integrantes del equipo terminar o bien realizar los \u00faltimos cambios sobre para Trabajos pr\u00e1cticos evaluables que se hayan realizado durante la \u00faltima semana.\",\n  { matchCase: true }\n);\naprovecharanHits.load(\"text\");\nawait context.sync();\nif (aprovecharanHits.items.length > 0) {\n  aprovecharanHits.items[0].insertText(\n    \"Se aprovechar\u00e1n los d\u00edas h\u00e1biles para que los integrantes del equipo terminen de  realizar los \u00faltimos cambios sobre los Trabajos pr\u00e1cticos evaluables que se hayan realizado durante la \u00faltima semana.\",\n    \"Replace\"\n  );\n  await context.sync();\n}\n\n// --- Change 3: \"hallan faltado\" -> \"hayan faltado\" ---\nconst hallanHits = body.search(\"hallan faltado\", { matchCase: true });\nhallanHits.load(\"text\");\nawait context.sync();\nif (hallanHits.items.length > 0) {\n  hallanHits.items[0].insertText(\"hayan faltado\", \"Replace\");\n  await context.sync();\n}\n", "ps1": "# Corecion TP4 (Gestion Linea Base)\n# 1) Insert a missing space after \"...fechas de ex\u00e1menes\" (before the\n#    existing spacing-adjusted run that precedes \"parciales\").\n# 2) Fix grammar in \"...terminar o bien realizar los \u00faltimos cambios\n#    sobre para Trabajos...\" -> \"...terminen de  realizar los \u00faltimos\n#    cambios sobre los Trabajos...\".\n# 3) Fix \"hallan faltado\" -> \"hayan faltado\".\n\n$d = $word.ActiveDocument\n\n# --- Change 1: add a space after \"... fechas de ex\u00e1menes\" ---\n$rng1 = $d.Content\n$rng1.Find.ClearFormatting()\n$rng1.Find.Text = \"h\u00e1bil posterior a las fechas de ex\u00e1menes\"\n$rng1.Find.Execute() | Out-Null\nif ($rng1.Find.Found) {\n    $rng1.Collapse(0)\n    $rng1.InsertAfter(\" \")\n}\n\n# --- Change 2: rewrite the \"Se aprovechar\u00e1n...\" sentence ---\n$rng2 = $d.Content\n$rng2.Find.ClearFormatting()\n$old2 = \"Se aprovechar\u00e1n los d\u00edas h\u00e1biles para que los integrantes del equipo terminar o bien realizar los \u00faltimos cambios sobre para Trabajos pr\u00e1cticos evaluables que se hayan realizado durante la \u00faltima semana.\"\n$new2 = \"Se aprovechar\u00e1n los d\u00edas h\u00e1biles para que los integrantes del equipo terminen de  realizar los \u00faltimos cambios sobre los Trabajos pr\u00e1cticos evaluables que se hayan realizado durante la \u00faltima semana.\"\n$rng2.Find.Text = $old2\n$rng2.Find.Replacement.Text = $new2\n$rng2.Find.Execute($old2, $false, $false, $false, $false, $false, $true, 1, $false, $new2, 2) | Out-Null\n\n# --- Change 3: \"hallan faltado\" -> \"hayan faltado\" ---\n$rng3 = $d.Content\n$rng3.Find.ClearFormatting()\n$old3 = \"hallan faltado\"\n$new3 = \"hayan faltado\"\n$rng3.Find.Text = $old3\n$rng3.Find.Replacement.Text = $new3\n$rng3.Find.Execute($old3, $false, $false, $false, $false, $false, $true, 1, $false, $new3, 2) | Out-Null\n"}
